# Regenerate save_data to use K (strikeouts) instead of Strike# column.
# Update column G ("K") values for each game row with the recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 3
    14 = 0
    15 = 2
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 3
    27 = 1
    29 = 2
    30 = 1
    32 = 1
    33 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
